# This workbook contains one week's worth of daily price observations per
# row, for "Vega Monumental Concepción - Coliflor". A new weekly record is
# inserted at row 116 (pushing the existing rows 116:220 down to 117:221),
# and the new row is populated with this week's values. Since every row
# shares the same Mercado/Región/Categoría/Unidad/Clasificación metadata,
# those columns are simply copied from the row that used to be row 116
# (now row 117) to keep them consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 116; this shifts old rows 116:220 down to
# 117:221 (so the old last row, 220, becomes the new row 221 automatically).
$ws.Rows.Item(116).Insert()

# Copy the constant metadata columns from row 117 (the row that used to be
# row 116) into the newly inserted row 116.
$ws.Range("A116").Value = $ws.Range("A117").Value2
$ws.Range("B116").Value = $ws.Range("B117").Value2
$ws.Range("C116").Value = $ws.Range("C117").Value2
$ws.Range("E116").Value = $ws.Range("E117").Value2
$ws.Range("F116").Value = $ws.Range("F117").Value2
$ws.Range("G116").Value = $ws.Range("G117").Value2
$ws.Range("H116").Value = $ws.Range("H117").Value2
$ws.Range("N116").Value = $ws.Range("N117").Value2
$ws.Range("Q116").Value = $ws.Range("Q117").Value2
$ws.Range("R116").Value = $ws.Range("R117").Value2

# New observation's own data.
$ws.Range("D116").Value = 44589
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 1000
$ws.Range("L116").Value = 1100
$ws.Range("M116").Value = 1038
$ws.Range("O116").Value = "Región Metropolitana"
$ws.Range("P116").Value = 1038
